$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "47.080.94"
$ws.Range("E2").Value2 = "  +1.18%  "

$ws.Range("D3").Value2 = "2.487.45"
$ws.Range("E3").Value2 = "  +0.86%  "

$ws.Range("E4").Value2 = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "321.66"
$ws.Range("E5").Value2 = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "108.58"
$ws.Range("E6").Value2 = "  +3.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.522"
$ws.Range("E7").Value2 = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.999"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.534"
$ws.Range("E9").Value2 = "  -0.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "38.79"
$ws.Range("E10").Value2 = "  +6.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0810"
$ws.Range("E11").Value2 = "  -0.73%  "

$ws.Range("E12").Value2 = "  +0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "18.20"
$ws.Range("E13").Value2 = "  -0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "7.13"
$ws.Range("E14").Value2 = "  +0.46%  "

$ws.Range("D15").Value2 = "2.876.49"
$ws.Range("E15").Value2 = "  +0.85%  "

$ws.Range("D16").Value2 = "2.488.16"
$ws.Range("E16").Value2 = "  +1.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.845"
$ws.Range("E17").Value2 = "  -0.07%  "

$ws.Range("D18").Value2 = "47.009.17"
$ws.Range("E18").Value2 = "  +1.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "12.69"
$ws.Range("E19").Value2 = "  -0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "6.60"
$ws.Range("E20").Value2 = "  +2.15%  "

$ws.Range("D21").Value2 = "0.0₃0933"
$ws.Range("E21").Value2 = "  -0.42%  "

$ws.Range("E22").Value2 = "  +15.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "70.55"
$ws.Range("E23").Value2 = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "245.62"
$ws.Range("E24").Value2 = "  -1.20%  "

$ws.Range("E25").Value2 = "  +0.91%  "

$ws.Range("E26").Value2 = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "25.68"
$ws.Range("E27").Value2 = "  -1.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "2.29"
$ws.Range("E28").Value2 = "  -1.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "10.05"
$ws.Range("E29").Value2 = "  +2.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.140"
$ws.Range("E30").Value2 = "  +8.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "35.04"
$ws.Range("E31").Value2 = "  -0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "49.90"
$ws.Range("E32").Value2 = "  +0.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "20.00"
$ws.Range("E33").Value2 = "  +2.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "5.39"
$ws.Range("E34").Value2 = "  +1.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.0782"
$ws.Range("E35").Value2 = "  +1.74%  "

$ws.Range("E36").Value2 = "  +0.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "4.67"
$ws.Range("E37").Value2 = "  +0.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "1.95"
$ws.Range("E38").Value2 = "  +2.19%  "

$ws.Range("E39").Value2 = "  -0.04%  "

$ws.Range("E40").Value2 = "  +0.07%  "

$ws.Range("E41").Value2 = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "118.49"
$ws.Range("E42").Value2 = "  -3.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "21.33"
$ws.Range("E43").Value2 = "  +3.08%  "

$ws.Range("E44").Value2 = "  +0.20%  "

$ws.Range("D45").Value2 = "1.979.91"
$ws.Range("E45").Value2 = "  -0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3.02"
$ws.Range("E46").Value2 = "  +1.16%  "

$ws.Range("E47").Value2 = "  -2.39%  "

$ws.Range("B48").Value2 = "FraxShare"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "9.04"
$ws.Range("E48").Value2 = "  +0.27%  "

$ws.Range("B49").Value2 = "Stacks"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.77"
$ws.Range("E49").Value2 = "  -1.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "5.13"
$ws.Range("E50").Value2 = "  -3.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "57.31"
$ws.Range("E51").Value2 = "  +4.57%  "
